$d = $word.ActiveDocument

# The document starts with 6 paragraphs:
#   1: "1) Hiện tại : " + "CARD_MONEY" + " >= " + "RECHARGE_MONEY"
#   2: "-> 2 số này có = nhau hay ko? Nếu không = nhau thì phải thêm 1 cột : " + ... + ". Khi đó:"
#   3: "CARD_MONEY" + " >= " + "RECHARGE_" + "REAL_" + "MONEY"
#   4: "CARD_MONEY" + " = " + "RECHARGE_MONEY"
#   5: "2) Hiện tại: Giá hiện tại cập nhật mỗi lần login -> có nên cập nhật mỗi lần quét thẻ ko?"
#   6: "3)Số lít xăng khi mua tính theo giá thực tế hay giá trên thẻ để trừ vào tổng kho" + "?" + " Hiện tạ" + "i tính theo giá trên thẻ"
#
# Target is 2 paragraphs:
#   1: "1) Hiện tại: Giá hiện tại cập nhật mỗi lần login -> có nên cập nhật mỗi lần quét thẻ ko?"
#      (as two runs: "1" and ") Hiện tại: ... ko?")
#   2: "2) Khi bán xăng ở trạm xăng: có nên hiển thị tổng số lít còn lại trong kho không?"

# --- Delete paragraphs 2 through 5 entirely: the leftover scratch-work
#     paragraphs plus the original question 2. This merges what is left of
#     paragraph 1 directly with old paragraph 6. ---
$p1 = $d.Paragraphs.Item(1)
$p5 = $d.Paragraphs.Item(5)
$d.Range($p1.Range.End, $p5.Range.End).Delete()

# --- Paragraph 1 now reads "1) Hiện tại : CARD_MONEY >= RECHARGE_MONEY".
#     Replace its text, keeping "1" and the rest of the new question as two
#     separate insertions so the leading "1" stays its own run. ---
$p1 = $d.Paragraphs.Item(1)
$rWhole = $d.Range($p1.Range.Start, $p1.Range.End)
$rWhole.Text = "1"

$p1 = $d.Paragraphs.Item(1)
$rInsert = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$rInsert.InsertAfter(") Hiện tại: Giá hiện tại cập nhật mỗi lần login -> có nên cập nhật mỗi lần quét thẻ ko?")

# --- Paragraph 2 (originally question 3) becomes the new question 2. ---
$p2 = $d.Paragraphs.Item(2)
$rQ2 = $d.Range($p2.Range.Start, $p2.Range.End)
$rQ2.Text = "2) Khi bán xăng ở trạm xăng: có nên hiển thị tổng số lít còn lại trong kho không?"
